$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "31.176.54", "  +2.01%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.940.39", "  +0.90%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.000", "  -0.45%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "242.23", "  +1.18%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.000", "  -0.38%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4798", "  -0.14%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2914", "  +1.18%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06788", "  +1.10%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "20.17", "  +7.24%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "104.30", "  -0.04%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07840", "  +1.13%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.941.39", "  +0.97%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.306", "  +1.16%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6988", "  +2.36%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "296.32", "  +11.20%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "31.152.82", "  +1.84%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "13.02", "  +2.21%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000007608", "  +0.98%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.575", "  +2.50%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.190.36", "  +0.85%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.000", "  -0.41%  "),
    @("BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.000", "  -0.39%  "),
    @("BitDAO", "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit", "0.4605", "  +0.69%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "6.417", "  +1.14%  "),
    @("Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "9.564", "  -0.76%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "169.07", "  +3.39%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "19.84", "  +4.19%  "),
    @("LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.097", "  -0.33%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.393", "  +0.31%  "),
    @("Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1009", "  -1.25%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.625", "  +0.58%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.536", "  +1.05%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.344", "  +1.72%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04847", "  +1.88%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7388", "  -0.23%  "),
    @("ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.134", "  +0.84%  "),
    @("HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.724", "  +1.54%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01965", "  +1.17%  "),
    @("FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.825", "  +7.45%  "),
    @("MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.634", "  -0.33%  "),
    @("Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "76.42", "  +1.03%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.041", "  +1.90%  "),
    @("TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.8733", "  +1.31%  "),
    @("TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4371", "  +1.79%  "),
    @("Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "105.79", "  -0.47%  "),
    @("PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.000", "  -0.39%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.032.93", "  +3.82%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.588", "  +0.19%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.231", "  +2.35%  "),
    @("Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1209", "  +0.52%  ")
)

$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
}

$wb.Save()